$wb = $excel.ActiveWorkbook

# ---- Sheet: BFS ----
$ws = $wb.Worksheets.Item("BFS")
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1,5).Value = "Distancia"
$ws.Cells.Item(2,2).Value = 19
$ws.Cells.Item(2,3).Value = 19
$ws.Cells.Item(2,4).Value = "[19]"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0.000213623046875
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(3,2).Value = 6
$ws.Cells.Item(3,3).Value = 13
$ws.Cells.Item(3,4).Value = "[6, 9, 13]"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 11
$ws.Cells.Item(3,7).Value = 1.363636363636364
$ws.Cells.Item(3,8).Value = 0.0001287460327148438
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(4,2).Value = 17
$ws.Cells.Item(4,3).Value = 31
$ws.Cells.Item(4,4).Value = "[17, 20, 31]"
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 12
$ws.Cells.Item(4,7).Value = 1.583333333333333
$ws.Cells.Item(4,8).Value = 0.0001022815704345703
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(5,2).Value = 14
$ws.Cells.Item(5,3).Value = 14
$ws.Cells.Item(5,4).Value = "[14]"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0.00007867813110351562
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = 12
$ws.Cells.Item(6,4).Value = "[3, 2, 5, 8, 12]"
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 15
$ws.Cells.Item(6,7).Value = 1.133333333333333
$ws.Cells.Item(6,8).Value = 0.00008821487426757812
$ws.Cells.Item(6,9).Value = 0

# ---- Sheet: DFS ----
$ws = $wb.Worksheets.Item("DFS")
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1,5).Value = "Distancia"
$ws.Cells.Item(2,2).Value = 19
$ws.Cells.Item(2,3).Value = 19
$ws.Cells.Item(2,4).Value = "[19]"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0.0000820159912109375
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(3,2).Value = 6
$ws.Cells.Item(3,3).Value = 13
$ws.Cells.Item(3,4).Value = "[6, 3, 2, 5, 4, 7, 8, 9, 28, 29, 13]"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 15
$ws.Cells.Item(3,7).Value = 0.1333333333333333
$ws.Cells.Item(3,8).Value = 0.00007796287536621094
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(4,2).Value = 17
$ws.Cells.Item(4,3).Value = 31
$ws.Cells.Item(4,4).Value = "[17, 13, 9, 6, 3, 2, 5, 4, 7, 8, 12, 11, 10, 14, 15, 16, 19, 18, 21, 22, 23, 24, 20, 31]"
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 25
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(4,8).Value = 0.00008344650268554688
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(5,2).Value = 14
$ws.Cells.Item(5,3).Value = 14
$ws.Cells.Item(5,4).Value = "[14]"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0.00007128715515136719
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = 12
$ws.Cells.Item(6,4).Value = "[3, 2, 5, 4, 7, 8, 9, 6, 27, 28, 29, 13, 12]"
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 16
$ws.Cells.Item(6,7).Value = 0.0625
$ws.Cells.Item(6,8).Value = 0.00007414817810058594
$ws.Cells.Item(6,9).Value = 0

# ---- Sheet: BCU ----
$ws = $wb.Worksheets.Item("BCU")
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1,5).Value = "Distancia"
$ws.Cells.Item(2,2).Value = 19
$ws.Cells.Item(2,3).Value = 19
$ws.Cells.Item(2,4).Value = "[19]"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 0.0001046657562255859
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(3,2).Value = 6
$ws.Cells.Item(3,3).Value = 13
$ws.Cells.Item(3,4).Value = "[6, 9, 13]"
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 9
$ws.Cells.Item(3,7).Value = 0.5714285714285714
$ws.Cells.Item(3,8).Value = 0.0002241134643554688
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(4,2).Value = 17
$ws.Cells.Item(4,3).Value = 31
$ws.Cells.Item(4,4).Value = "[17, 30, 31]"
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 7
$ws.Cells.Item(4,7).Value = 0.4
$ws.Cells.Item(4,8).Value = 0.0001020431518554688
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(5,2).Value = 14
$ws.Cells.Item(5,3).Value = 14
$ws.Cells.Item(5,4).Value = "[14]"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0
$ws.Cells.Item(5,8).Value = 0.0001242160797119141
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = 12
$ws.Cells.Item(6,4).Value = "[3, 6, 5, 8, 12]"
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 18
$ws.Cells.Item(6,7).Value = 0.7727272727272727
$ws.Cells.Item(6,8).Value = 0.0001125335693359375
$ws.Cells.Item(6,9).Value = 0

# ---- Sheet: A_Estrela_Euclidiano ----
$ws = $wb.Worksheets.Item("A_Estrela_Euclidiano")
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1,5).Value = "Distancia"
$ws.Cells.Item(2,2).Value = 19
$ws.Cells.Item(2,3).Value = 19
$ws.Cells.Item(2,4).Value = "[19]"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.0625
$ws.Cells.Item(2,8).Value = 0.0001130104064941406
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(3,2).Value = 6
$ws.Cells.Item(3,3).Value = 13
$ws.Cells.Item(3,4).Value = "[6, 9, 13]"
$ws.Cells.Item(3,5).Value = 164
$ws.Cells.Item(3,6).Value = 9
$ws.Cells.Item(3,7).Value = 3.0625
$ws.Cells.Item(3,8).Value = 0.0001275539398193359
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(4,2).Value = 17
$ws.Cells.Item(4,3).Value = 31
$ws.Cells.Item(4,4).Value = "[17, 30, 31]"
$ws.Cells.Item(4,5).Value = 166
$ws.Cells.Item(4,6).Value = 6
$ws.Cells.Item(4,7).Value = 3.0625
$ws.Cells.Item(4,8).Value = 0.0001168251037597656
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(5,2).Value = 14
$ws.Cells.Item(5,3).Value = 14
$ws.Cells.Item(5,4).Value = "[14]"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.0625
$ws.Cells.Item(5,8).Value = 0.0001053810119628906
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = 12
$ws.Cells.Item(6,4).Value = "[3, 6, 5, 8, 12]"
$ws.Cells.Item(6,5).Value = 349
$ws.Cells.Item(6,6).Value = 16
$ws.Cells.Item(6,7).Value = 3.0625
$ws.Cells.Item(6,8).Value = 0.0004928112030029297
$ws.Cells.Item(6,9).Value = 0

# ---- Sheet: A_Estrela_Haversiano ----
$ws = $wb.Worksheets.Item("A_Estrela_Haversiano")
$ws.Columns.Item(5).Insert()
$ws.Cells.Item(1,5).Value = "Distancia"
$ws.Cells.Item(2,2).Value = 19
$ws.Cells.Item(2,3).Value = 19
$ws.Cells.Item(2,4).Value = "[19]"
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.0625
$ws.Cells.Item(2,8).Value = 0.003904104232788086
$ws.Cells.Item(2,9).Value = 0
$ws.Cells.Item(3,2).Value = 6
$ws.Cells.Item(3,3).Value = 13
$ws.Cells.Item(3,4).Value = "[6, 9, 13]"
$ws.Cells.Item(3,5).Value = 164
$ws.Cells.Item(3,6).Value = 9
$ws.Cells.Item(3,7).Value = 3.0625
$ws.Cells.Item(3,8).Value = 0.0001533031463623047
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(4,2).Value = 17
$ws.Cells.Item(4,3).Value = 31
$ws.Cells.Item(4,4).Value = "[17, 30, 31]"
$ws.Cells.Item(4,5).Value = 166
$ws.Cells.Item(4,6).Value = 6
$ws.Cells.Item(4,7).Value = 3.0625
$ws.Cells.Item(4,8).Value = 0.0001726150512695312
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(5,2).Value = 14
$ws.Cells.Item(5,3).Value = 14
$ws.Cells.Item(5,4).Value = "[14]"
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.0625
$ws.Cells.Item(5,8).Value = 0.0001351833343505859
$ws.Cells.Item(5,9).Value = 0
$ws.Cells.Item(6,2).Value = 3
$ws.Cells.Item(6,3).Value = 12
$ws.Cells.Item(6,4).Value = "[3, 6, 5, 8, 12]"
$ws.Cells.Item(6,5).Value = 349
$ws.Cells.Item(6,6).Value = 16
$ws.Cells.Item(6,7).Value = 3.0625
$ws.Cells.Item(6,8).Value = 0.0001521110534667969
$ws.Cells.Item(6,9).Value = 0

